# Update values in column C (RandomForest imputed results) for specific rows
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3   = -12.4664
    14  = -13.6226
    21  = -13.01700000000001
    23  = -12.1756
    25  = -11.3914
    26  = -12.17830000000001
    29  = -11.7003
    53  = -10.86040000000001
    57  = -14.26109999999999
    59  = -12.5142
    69  = -10.89559999999999
    79  = -11.52330000000001
    83  = -13.5781
    91  = -12.47270000000001
    93  = -10.347
    103 = -12.36949999999999
}

foreach ($row in $updates.Keys) {
    $ws.Range("C$row").Value2 = $updates[$row]
}
